$wb = $excel.ActiveWorkbook
$wsBudget = $wb.Worksheets.Item("budget")
$wsExpenses = $wb.Worksheets.Item("expenses")

# ---------------------------------------------------------------------------
# "budget" sheet: update existing row 2, and append new rows 3-13
# ---------------------------------------------------------------------------

# Row 2 values changed (A and B only; rest stay the same)
$wsBudget.Cells.Item(2, 1).Value = 35000
$wsBudget.Cells.Item(2, 2).Value = 15000

$budgetRows = @(
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 3
    @(20000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 190000), # row 4
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 5
    @(15000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 185000), # row 6
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 7
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 8
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 9
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 10
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 11
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000), # row 12
    @(10000, 20000, 12000, 5000, 10000, 8000, 15000, 10000, 45000, 30000, 15000, 180000)  # row 13
)

for ($i = 0; $i -lt $budgetRows.Length; $i++) {
    $rowNum = $i + 3
    $rowValues = $budgetRows[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $wsBudget.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

# Move the selection on the "budget" sheet to L21, matching the saved view
$wsBudget.Range("L21").Select()

# ---------------------------------------------------------------------------
# "expenses" sheet: a handful of cells (loans column + totals) were revised
# ---------------------------------------------------------------------------

$wsExpenses.Cells.Item(2, 1).Value = 12000
$wsExpenses.Cells.Item(2, 12).Value = 128000

$wsExpenses.Cells.Item(3, 1).Value = 15000
$wsExpenses.Cells.Item(3, 12).Value = 151000

$wsExpenses.Cells.Item(4, 1).Value = 15000
$wsExpenses.Cells.Item(4, 12).Value = 147600

$wsExpenses.Cells.Item(6, 1).Value = 15000
$wsExpenses.Cells.Item(6, 12).Value = 132550

$wsExpenses.Cells.Item(10, 1).Value = 12000
$wsExpenses.Cells.Item(10, 12).Value = 158821
